$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8356
$ws1.Range("F5").Value = 8356
$ws1.Range("F7").Value = 86
$ws1.Range("F8").Value = 408
$ws1.Range("F9").Value = 7322
$ws1.Range("F10").Value = 584
$ws1.Range("F11").Value = 506
$ws1.Range("F15").Value = 233
$ws1.Range("F18").Value = 142
$ws1.Range("F19").Value = 12118
$ws1.Range("F22").Value = 2449
$ws1.Range("F23").Value = 3491
$ws1.Range("F26").Value = 2894
$ws1.Range("F27").Value = 108
$ws1.Range("F29").Value = 37
$ws1.Range("F30").Value = 3344
$ws1.Range("F32").Value = 340
$ws1.Range("F33").Value = 1706
$ws1.Range("F35").Value = 127
$ws1.Range("F36").Value = 5997
$ws1.Range("F37").Value = 97
$ws1.Range("F38").Value = 1829
$ws1.Range("F39").Value = 1253
$ws1.Range("F40").Value = 30
$ws1.Range("F41").Value = 895
$ws1.Range("F45").Value = 196
$ws1.Range("F47").Value = 1113
$ws1.Range("F48").Value = 1576
$ws1.Range("F50").Value = 114

# Sheet "演出" (index 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 7
$ws2.Range("F22").Value = 72
$ws2.Range("F25").Value = 4

# Sheet "本地生活" (index 3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 316
$ws3.Range("F3").Value = 460
$ws3.Range("F4").Value = 12

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 316
$ws4.Range("F6").Value = 7
$ws4.Range("F8").Value = 8356
$ws4.Range("F9").Value = 86
$ws4.Range("F11").Value = 408
$ws4.Range("F12").Value = 7322
$ws4.Range("F13").Value = 7322
$ws4.Range("F14").Value = 584
$ws4.Range("F15").Value = 506
$ws4.Range("F18").Value = 233
$ws4.Range("F21").Value = 142
$ws4.Range("F23").Value = 12118
$ws4.Range("F27").Value = 2449
$ws4.Range("F28").Value = 2449
$ws4.Range("F29").Value = 3491
$ws4.Range("F30").Value = 108
$ws4.Range("F32").Value = 37
$ws4.Range("F34").Value = 3344
$ws4.Range("F36").Value = 340
$ws4.Range("F37").Value = 1706
$ws4.Range("F39").Value = 127
$ws4.Range("F40").Value = 5997
$ws4.Range("F41").Value = 72
$ws4.Range("F42").Value = 1829
$ws4.Range("F44").Value = 1253
$ws4.Range("F45").Value = 30
$ws4.Range("F46").Value = 895
$ws4.Range("F48").Value = 196
$ws4.Range("F50").Value = 1113
$ws4.Range("F51").Value = 1576
$ws4.Range("F52").Value = 114
